$d = $word.ActiveDocument

# Locate the "Jeremy Schroeder" text in the author list and replace it with
# "Alex Hall", preserving the original run layout (separate runs for the
# first name, the space, and the last name) by inserting OOXML markup
# instead of doing a plain text replace (which would collapse the runs).
$text = $d.Content.Text
$oldName = "Jeremy Schroeder"
$idx = $text.IndexOf($oldName)

if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + $oldName.Length)

    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">Alex</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Hall</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

    $r.InsertXML($xml)
}
